$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Values for the new column J (rows 2-16)
$jValues = @(111, 112, 113, 121, 122, 123, 131, 132, 133, 1, 2, 3, 1, 2, 3)

for ($i = 0; $i -lt $jValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}

# Set custom column widths for columns C through G
# (input values chosen so the engine's internal character-width
# quantization lands as close as possible to the target stored widths)
$ws.Columns.Item(3).ColumnWidth = 14.285714285714286
$ws.Columns.Item(4).ColumnWidth = 11.857142857142858
$ws.Columns.Item(5).ColumnWidth = 16.57142857142857
$ws.Columns.Item(6).ColumnWidth = 22.0
$ws.Columns.Item(7).ColumnWidth = 11.0

# Update the active selection to E16
$ws.Range("E16").Select()
